# BubbleSort.xlsx update — refreshed timing measurements for the
# Bubble Sort benchmark, plus a couple of chart touch-ups that came
# along with the re-save (trendline/series line weights, hiding the
# R-squared trendline label).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Updated benchmark timings (column B) ------------------------------
$ws.Range("B1").Value  = 588
$ws.Range("B2").Value  = 2130
$ws.Range("B3").Value  = 4707
$ws.Range("B4").Value  = 12610
$ws.Range("B5").Value  = 14570
$ws.Range("B6").Value  = 21869
$ws.Range("B7").Value  = 34813
$ws.Range("B8").Value  = 43468
$ws.Range("B9").Value  = 56638
$ws.Range("B10").Value = 69338

# --- Chart touch-ups -----------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)

# Thicker plotted line for the data series.
$series.Border.Weight = 2

# Thicker trendline, and stop displaying the R-squared value on it.
$trendline = $series.Trendlines().Item(1)
$trendline.Border.Weight = 3
$trendline.DisplayRSquared = $false

# --- Restore the active selection used when the file was last saved ------
[void]$ws.Range("C17").Select()
